# Revert commit 6dc3ea4: restore original Abstract/Authors text for rows 3, 13, 14
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s_D3 = @'
Background
id="Par1">The coronavirus disease 2019 (COVID-19) pandemic caused by severe acute respiratory syndrome coronavirus 2 (SARS-CoV-2) infection has profoundly affected the lives of millions of people.

 To date, there is no approved vaccine or specific drug to prevent or treat COVID-19, while the infection is globally spreading at an alarming rate.

 Because the development of effective vaccines or novel drugs could take several months (if not years), repurposing existing drugs is considered a more efficient strategy that could save lives now.

 Statins constitute a class of lipid-lowering drugs with proven safety profiles and various known beneficial pleiotropic effects.

 Our previous investigations showed that statins have antiviral effects and are involved in the process of wound healing in the lung.

 This triggered us to evaluate if statin use reduces mortality in COVID-19 patients.


Results
id="Par2">After initial recruitment of 459 patients with COVID-19 (Shiraz province, Iran) and careful consideration of the exclusion criteria, a total of 150 patients, of which 75 received statins, were included in our retrospective study.

 Cox proportional-hazards regression models were used to estimate the association between statin use and rate of death.

 After propensity score matching, we found that statin use appeared to be associated with a lower risk of morbidity [HR = 0.85, 95% CI = (0.02, 3.93), P = 0.762] and lower risk of death [(HR = 0.76; 95% CI = (0.16, 3.72), P = 0.735)]; however, these associations did not reach statistical significance.

 Furthermore, statin use reduced the chance of being subjected to mechanical ventilation [OR = 0.96, 95% CI = (0.61–2.99), P = 0.942] and patients on statins showed a more normal computed tomography (CT) scan result [OR = 0.41, 95% CI = (0.07–2.33), P = 0.312].


Conclusions
id="Par3">Although we could not demonstrate a significant association between statin use and a reduction in mortality in patients with COVID19, we do feel that our results are promising and of clinical relevance and warrant the need for prospective randomized controlled trials and extensive retrospective studies to further evaluate and validate the potential beneficial effects of statin treatment on clinical symptoms and mortality rates associated with COVID-19.

'@
$ws.Range("D3").Value2 = $s_D3

$s_E3 = @'
[Payam%Peymani%NULL%1,    Tania%Dehesh%NULL%1,    Farnaz%Aligolighasemabadi%NULL%1,    Mohammadamin%Sadeghdoust%NULL%1,    Katarzyna%Kotfis%NULL%1,    Mazaher%Ahmadi%NULL%1,    Parvaneh%Mehrbod%NULL%1,    Pooya%Iranpour%NULL%1,    Sanaz%Dastghaib%NULL%1,    Ahmad%Nasimian%NULL%1,    Amir%Ravandi%NULL%1,    Biniam%Kidane%NULL%1,    Naseer%Ahmed%NULL%1,    Pawan%Sharma%NULL%1,    Shahla%Shojaei%NULL%1,    Kamran%Bagheri Lankarani%NULL%1,    Andrzej%Madej%NULL%1,    Nima%Rezaei%NULL%1,    Tayyebeh%Madrakian%NULL%1,    Marek J.%Los%NULL%1,    Hagar Ibrahim%Labouta%NULL%1,    Pooneh%Mokarram%NULL%1,    Saeid%Ghavami%saeid.ghavami@umanitoba.ca%1]
'@
$ws.Range("E3").Value2 = $s_E3

$s_D13 = @'
id="Par1">The coronavirus disease 2019 (COVID-19) can result in a hyperinflammatory state, leading to acute respiratory distress syndrome (ARDS), myocardial injury, and thrombotic complications, among other sequelae.

 Statins, which are known to have anti-inflammatory and antithrombotic properties, have been studied in the setting of other viral infections, but their benefit has not been assessed in COVID-19. This is a retrospective analysis of patients admitted with COVID-19 from February 1st through May 12th, 2020 with study period ending on June 11th, 2020. Antecedent statin use was assessed using medication information available in the electronic medical record.

 We constructed a multivariable logistic regression model to predict the propensity of receiving statins, adjusting for baseline sociodemographic and clinical characteristics, and outpatient medications.

 The primary endpoint includes in-hospital mortality within 30 days.

 A total of 2626 patients were admitted during the study period, of whom 951 (36.2%) were antecedent statin users.

 Among 1296 patients (648 statin users, 648 non-statin users) identified with 1:1 propensity-score matching, statin use is significantly associated with lower odds of the primary endpoint in the propensity-matched cohort (OR 0.47, 95% CI 0.36–0.62, p &lt; 0.001).

 We conclude that antecedent statin use in patients hospitalized with COVID-19 is associated with lower inpatient mortality.


'@
$ws.Range("D13").Value2 = $s_D13

$s_E13 = @'
[Aakriti%Gupta%ag3786@cumc.columbia.edu%1,    Mahesh V.%Madhavan%NULL%2,    Mahesh V.%Madhavan%NULL%0,    Timothy J.%Poterucha%NULL%2,    Timothy J.%Poterucha%NULL%0,    Ersilia M.%DeFilippis%NULL%1,    Jessica A.%Hennessey%NULL%1,    Bjorn%Redfors%NULL%1,    Christina%Eckhardt%NULL%1,    Behnood%Bikdeli%NULL%1,    Jonathan%Platt%NULL%1,    Ani%Nalbandian%NULL%1,    Pierre%Elias%NULL%2,    Pierre%Elias%NULL%0,    Matthew J.%Cummings%NULL%2,    Matthew J.%Cummings%NULL%0,    Shayan N.%Nouri%NULL%1,    Matthew%Lawlor%NULL%1,    Lauren S.%Ranard%NULL%1,    Jianhua%Li%NULL%2,    Jianhua%Li%NULL%0,    Claudia%Boyle%NULL%1,    Raymond%Givens%NULL%1,    Daniel%Brodie%NULL%0,    Harlan M.%Krumholz%NULL%2,    Harlan M.%Krumholz%NULL%0,    Gregg W.%Stone%NULL%2,    Gregg W.%Stone%NULL%0,    Sanjum S.%Sethi%NULL%1,    Daniel%Burkhoff%NULL%1,    Nir%Uriel%NULL%1,    Allan%Schwartz%NULL%0,    Martin B.%Leon%NULL%1,    Ajay J.%Kirtane%NULL%1,    Elaine Y.%Wan%eyw2003@cumc.columbia.edu%1,    Sahil A.%Parikh%NULL%2,    Sahil A.%Parikh%NULL%0]
'@
$ws.Range("E13").Value2 = $s_E13

$s_D14 = @'
Aims/hypothesis
id="Par1">This is an update of the results from the previous report of the CORONADO (Coronavirus SARS-CoV-2 and Diabetes Outcomes) study, which aims to describe the outcomes and prognostic factors in patients with diabetes hospitalised for coronavirus disease-2019 (COVID-19).


Methods
id="Par2">The CORONADO initiative is a French nationwide multicentre study of patients with diabetes hospitalised for COVID-19 with a 28-day follow-up.

 The patients were screened after hospital admission from 10 March to 10 April 2020. We mainly focused on hospital discharge and death within 28 days.


Results
id="Par3">We included 2796 participants: 63.7% men, mean age 69.7 ± 13.2 years, median BMI (25th–75th percentile) 28.4 (25.0–32.4) kg/m2.

 Microvascular and macrovascular diabetic complications were found in 44.2% and 38.6% of participants, respectively.

 Within 28 days, 1404 (50.2%; 95% CI 48.3%, 52.1%) were discharged from hospital with a median duration of hospital stay of 9 (5–14) days, while 577 participants died (20.6%; 95% CI 19.2%, 22.2%).

 In multivariable models, younger age, routine metformin therapy and longer symptom duration on admission were positively associated with discharge.

 History of microvascular complications, anticoagulant routine therapy, dyspnoea on admission, and higher aspartate aminotransferase, white cell count and C-reactive protein levels were associated with a reduced chance of discharge.

 Factors associated with death within 28 days mirrored those associated with discharge, and also included routine treatment by insulin and statin as deleterious factors.


Conclusions/interpretation
id="Par4">In patients with diabetes hospitalised for COVID-19, we established prognostic factors for hospital discharge and death that could help clinicians in this pandemic period.


Trial registration
id="Par5">Clinicaltrials.

gov identifier: NCT04324736
Graphical abstract
id="Par6">
            
          
Supplementary Information
The online version contains peer-reviewed but unedited supplementary material available at 10.1007/s00125-020-05351-w.



'@
$ws.Range("D14").Value2 = $s_D14

$s_E14 = @'
[Matthieu%Wargny%NULL%1,    Louis%Potier%NULL%2,    Louis%Potier%NULL%0,    Pierre%Gourdy%NULL%2,    Pierre%Gourdy%NULL%0,    Matthieu%Pichelin%NULL%2,    Matthieu%Pichelin%NULL%0,    Coralie%Amadou%NULL%2,    Coralie%Amadou%NULL%0,    Pierre-Yves%Benhamou%NULL%2,    Pierre-Yves%Benhamou%NULL%0,    Jean-Baptiste%Bonnet%NULL%1,    Lyse%Bordier%NULL%1,    Olivier%Bourron%NULL%1,    Claude%Chaumeil%NULL%2,    Claude%Chaumeil%NULL%0,    Nicolas%Chevalier%NULL%1,    Patrice%Darmon%NULL%2,    Patrice%Darmon%NULL%0,    Blandine%Delenne%NULL%2,    Blandine%Delenne%NULL%0,    Delphine%Demarsy%NULL%1,    Marie%Dumas%NULL%1,    Olivier%Dupuy%NULL%1,    Anna%Flaus-Furmaniuk%NULL%2,    Anna%Flaus-Furmaniuk%NULL%0,    Jean-François%Gautier%NULL%1,    Anne-Marie%Guedj%NULL%2,    Anne-Marie%Guedj%NULL%0,    Nathalie%Jeandidier%NULL%1,    Etienne%Larger%NULL%1,    Jean-Philippe%Le Berre%NULL%1,    Myriam%Lungo%NULL%1,    Nathanaëlle%Montanier%NULL%1,    Philippe%Moulin%NULL%1,    Françoise%Plat%NULL%1,    Vincent%Rigalleau%NULL%1,    René%Robert%NULL%1,    Dominique%Seret-Bégué%NULL%2,    Dominique%Seret-Bégué%NULL%0,    Pierre%Sérusclat%NULL%1,    Sarra%Smati%NULL%1,    Jean-François%Thébaut%NULL%1,    Blandine%Tramunt%NULL%2,    Blandine%Tramunt%NULL%0,    Camille%Vatier%NULL%1,    Fritz-Line%Velayoudom%NULL%2,    Fritz-Line%Velayoudom%NULL%0,    Bruno%Vergès%NULL%1,    Patrice%Winiszewski%NULL%1,    Audrey%Zabulon%NULL%1,    Pierre-Antoine%Gourraud%NULL%1,    Ronan%Roussel%NULL%2,    Ronan%Roussel%NULL%0,    Bertrand%Cariou%bertrand.cariou@univ-nantes.fr%2,    Bertrand%Cariou%bertrand.cariou@univ-nantes.fr%0,    Samy%Hadjadj%samy.hadjadj@univ-nantes.fr%1,    NULL%NULL%NULL%0,    NULL%NULL%NULL%0]
'@
$ws.Range("E14").Value2 = $s_E14

